$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 188 (shifts 188..300 down to 189..301).
$ws.Rows.Item(188).Insert()

# The new row 188 re-uses the same "static" attributes (market/region/category/etc.)
# as the record that is now sitting at row 189, but carries its own date and
# volume/price observations.
$ws.Range("A188").Value = 4
$ws.Range("B188").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C188").Value = "Los Lagos"
$ws.Range("D188").Value = 44777
$ws.Range("E188").Value = 10
$ws.Range("F188").Value = 100112043
$ws.Range("G188").Value = "Pepino ensalada"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 200
$ws.Range("K188").Value = 25000
$ws.Range("L188").Value = 25000
$ws.Range("M188").Value = 25000
$ws.Range("N188").Value = "$/caja 60 unidades"
$ws.Range("O188").Value = "Región de Arica y Parinacota"
$ws.Range("P188").Value = 417
$ws.Range("Q188").Value = 60
$ws.Range("R188").Value = "Hortaliza"
